# Scheduled-runner refresh of Universalis market data for the "Faerie Profits"
# workbook. Updates the computed price/profit columns on each job sheet
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR):
#   H currentAveragePrice   I currentAveragePriceNQ   J currentAveragePriceHQ
#   K LevePriceNQ           L LevePriceHQ
#   M LeveProfitNQ          N LeveProfitHQ

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 89.63636
$ws.Range("I9").Value = 94.625
$ws.Range("J9").Value = 76.333336
$ws.Range("K9").Value = 94.625
$ws.Range("L9").Value = 76.333336
$ws.Range("M9").Value = 74.375
$ws.Range("N9").Value = -414.333336
$ws.Range("H92").Value = 1158
$ws.Range("I92").Value = 1000.05884
$ws.Range("K92").Value = 1000.05884
$ws.Range("M92").Value = 247.94116
$ws.Range("H135").Value = 7205.9565
$ws.Range("I135").Value = 6032.3335
$ws.Range("J135").Value = 11431
$ws.Range("K135").Value = 54291.0015
$ws.Range("L135").Value = 102879
$ws.Range("M135").Value = -51756.0015
$ws.Range("N135").Value = -107949
$ws.Range("H138").Value = 291157.12
$ws.Range("I138").Value = 33522.71
$ws.Range("J138").Value = 1432109.6
$ws.Range("K138").Value = 100568.13
$ws.Range("L138").Value = 4296328.800000001
$ws.Range("M138").Value = -95428.13
$ws.Range("N138").Value = -4306608.800000001

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2239.3125
$ws.Range("I45").Value = 1987.7858
$ws.Range("K45").Value = 1987.7858
$ws.Range("M45").Value = -1610.7858
$ws.Range("H63").Value = 2688.2222
$ws.Range("I63").Value = 2365.8333
$ws.Range("J63").Value = 3333
$ws.Range("K63").Value = 2365.8333
$ws.Range("L63").Value = 3333
$ws.Range("M63").Value = -1679.8333
$ws.Range("N63").Value = -4705
$ws.Range("H66").Value = 2688.2222
$ws.Range("I66").Value = 2365.8333
$ws.Range("J66").Value = 3333
$ws.Range("K66").Value = 11829.1665
$ws.Range("L66").Value = 16665
$ws.Range("M66").Value = -8397.166499999999
$ws.Range("N66").Value = -23529
$ws.Range("H97").Value = 2128.7368
$ws.Range("I97").Value = 1430.75
$ws.Range("J97").Value = 5851.3335
$ws.Range("K97").Value = 1430.75
$ws.Range("L97").Value = 5851.3335
$ws.Range("M97").Value = -934.75
$ws.Range("N97").Value = -6843.3335
$ws.Range("H122").Value = 1616.85
$ws.Range("I122").Value = 1491.421
$ws.Range("K122").Value = 4474.263
$ws.Range("M122").Value = -2024.263
$ws.Range("H132").Value = 1931.1666
$ws.Range("I132").Value = 1408.7273
$ws.Range("J132").Value = 7678
$ws.Range("K132").Value = 4226.1819
$ws.Range("L132").Value = 23034
$ws.Range("M132").Value = -1696.1819
$ws.Range("N132").Value = -28094

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4294.7144
$ws.Range("I94").Value = 3349.5
$ws.Range("K94").Value = 3349.5
$ws.Range("M94").Value = -2898.5
$ws.Range("H99").Value = 3059.5715
$ws.Range("I99").Value = 2241.1333
$ws.Range("K99").Value = 2241.1333
$ws.Range("M99").Value = -743.1333
$ws.Range("H105").Value = 4816.0557
$ws.Range("I105").Value = 3334.9285
$ws.Range("K105").Value = 3334.9285
$ws.Range("M105").Value = -1587.9285

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 39962.25
$ws.Range("I51").Value = 24999.5
$ws.Range("J51").Value = 54925
$ws.Range("K51").Value = 24999.5
$ws.Range("L51").Value = 54925
$ws.Range("M51").Value = -24263.5
$ws.Range("N51").Value = -56397
$ws.Range("H61").Value = 39962.25
$ws.Range("I61").Value = 24999.5
$ws.Range("J61").Value = 54925
$ws.Range("K61").Value = 24999.5
$ws.Range("L61").Value = 54925
$ws.Range("M61").Value = -24651.5
$ws.Range("N61").Value = -55621
$ws.Range("H134").Value = 2418.2766
$ws.Range("I134").Value = 1271.7142
$ws.Range("K134").Value = 3815.1426
$ws.Range("M134").Value = -1280.1426
$ws.Range("H141").Value = 121493.555
$ws.Range("J141").Value = 131680.25
$ws.Range("L141").Value = 131680.25
$ws.Range("N141").Value = -142040.25

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 142858430
$ws.Range("J23").Value = 166668140
$ws.Range("L23").Value = 500004420
$ws.Range("N23").Value = -500004890
$ws.Range("H38").Value = 366.23077
$ws.Range("I38").Value = 405.66666
$ws.Range("K38").Value = 1216.99998
$ws.Range("M38").Value = -869.9999800000001
$ws.Range("H69").Value = 4011.5
$ws.Range("I69").Value = 4011.5
$ws.Range("K69").Value = 12034.5
$ws.Range("M69").Value = -11223.5
$ws.Range("H72").Value = 4011.5
$ws.Range("I72").Value = 4011.5
$ws.Range("K72").Value = 36103.5
$ws.Range("M72").Value = -32047.5
$ws.Range("H120").Value = 4666.3335
$ws.Range("I120").Value = 4666.3335
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 13999.0005
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -9161.000499999998
$ws.Range("N120").ClearContents()
$ws.Range("H132").Value = 3477.4348
$ws.Range("I132").Value = 1133
$ws.Range("J132").Value = 3829.1
$ws.Range("K132").Value = 10197
$ws.Range("L132").Value = 34461.9
$ws.Range("M132").Value = -7667
$ws.Range("N132").Value = -39521.9

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 50587
$ws.Range("J123").Value = 50587
$ws.Range("L123").Value = 50587
$ws.Range("N123").Value = -55487
$ws.Range("H126").Value = 4864.5
$ws.Range("J126").Value = 5004.6665
$ws.Range("L126").Value = 15013.9995
$ws.Range("N126").Value = -19953.9995
$ws.Range("H132").Value = 7410177
$ws.Range("I132").Value = 8549771
$ws.Range("K132").Value = 25649313
$ws.Range("M132").Value = -25646783

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5723.3335
$ws.Range("I46").Value = 2362.1538
$ws.Range("J46").Value = 8293.647000000001
$ws.Range("K46").Value = 2362.1538
$ws.Range("L46").Value = 8293.647000000001
$ws.Range("M46").Value = -2174.1538
$ws.Range("N46").Value = -8669.647000000001
$ws.Range("H55").Value = 218.8
$ws.Range("I55").Value = 129.57143
$ws.Range("J55").Value = 296.875
$ws.Range("K55").Value = 129.57143
$ws.Range("L55").Value = 296.875
$ws.Range("M55").Value = 43.42857000000001
$ws.Range("N55").Value = -642.875
$ws.Range("H98").Value = 68499.664
$ws.Range("J98").Value = 68499.664
$ws.Range("L98").Value = 68499.664
$ws.Range("N98").Value = -74489.664
$ws.Range("H112").Value = 82218.5
$ws.Range("J112").Value = 82218.5
$ws.Range("L112").Value = 82218.5
$ws.Range("N112").Value = -85172.5
$ws.Range("H132").Value = 2951.7659
$ws.Range("I132").Value = 2904.077
$ws.Range("J132").Value = 3184.25
$ws.Range("K132").Value = 8712.231
$ws.Range("L132").Value = 9552.75
$ws.Range("M132").Value = -6182.231
$ws.Range("N132").Value = -14612.75
$ws.Range("H136").Value = 4119.8887
$ws.Range("I136").Value = 3722.8518
$ws.Range("K136").Value = 11168.5554
$ws.Range("M136").Value = -8618.555399999999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4568.5713
$ws.Range("I96").Value = 4245.25
$ws.Range("K96").Value = 4245.25
$ws.Range("M96").Value = -2872.25
$ws.Range("H100").Value = 1908.7241
$ws.Range("I100").Value = 1554.5652
$ws.Range("J100").Value = 3266.3333
$ws.Range("K100").Value = 3109.1304
$ws.Range("L100").Value = 6532.6666
$ws.Range("M100").Value = -2568.1304
$ws.Range("N100").Value = -7614.6666
$ws.Range("H136").Value = 15240.333
$ws.Range("I136").Value = 18583.834
$ws.Range("J136").Value = 1866.3334
$ws.Range("K136").Value = 55751.50199999999
$ws.Range("L136").Value = 5599.0002
$ws.Range("M136").Value = -53201.50199999999
$ws.Range("N136").Value = -10699.0002
